# Apply the edits described by the diff to Sheet1:
# - Populate two new header columns (N1/O1) with jsonpath-style column
#   names, and the matching example values in the first data row (N2/O2)
# - Grow the header row so the extra header text wraps/fits
# - Move the selection to the newly added O1 header cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cells - order matters for shared-string table index assignment,
# matching the order the values were originally typed in.
$ws.Range("O2").Value = "QnAYesNoBot"
$ws.Range("N2").Value = "Test"
$ws.Range("N1").Value = "clientFilterValues"
$ws.Range("O1").Value = "elicitResponse.responsebot_hook"

# Row 1 height change (34 -> 51) to accommodate the new headers
$ws.Rows.Item(1).RowHeight = 51

# Update the selection / view to the new O1 header cell
$ws.Range("O1").Select()
